$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing Text format so
# numeric-looking strings (e.g. "1.002", "0.02310") are preserved literally,
# matching the original inlineStr/text cells in the workbook.
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextCell 'D2' '27.969.60'
Set-TextCell 'E2' '  -1.78%  '
Set-TextCell 'D3' '1.763.43'
Set-TextCell 'E3' '  -3.76%  '
Set-TextCell 'D4' '1.002'
Set-TextCell 'E4' '  +0.88%  '
Set-TextCell 'D5' '320.49'
Set-TextCell 'E5' '  -2.12%  '
Set-TextCell 'D6' '1.001'
Set-TextCell 'E6' '  +0.87%  '
Set-TextCell 'D7' '0.4247'
Set-TextCell 'E7' '  -4.77%  '
Set-TextCell 'D8' '0.3609'
Set-TextCell 'E8' '  -4.88%  '
Set-TextCell 'D9' '43.73'
Set-TextCell 'E9' '  -3.86%  '
Set-TextCell 'D10' '0.07444'
Set-TextCell 'E10' '  -4.39%  '
Set-TextCell 'D11' '1.097'
Set-TextCell 'E11' '  -4.01%  '
Set-TextCell 'D12' '1.002'
Set-TextCell 'E12' '  +1.01%  '
Set-TextCell 'D13' '21.06'
Set-TextCell 'E13' '  -5.77%  '
Set-TextCell 'D14' '6.070'
Set-TextCell 'E14' '  -4.11%  '
Set-TextCell 'D15' '7.315'
Set-TextCell 'E15' '  -3.36%  '
Set-TextCell 'D16' '1.788.74'
Set-TextCell 'E16' '  -1.90%  '
Set-TextCell 'D18' '0.00001056'
Set-TextCell 'E18' '  -2.83%  '
Set-TextCell 'D19' '0.06385'
Set-TextCell 'E19' '  -0.06%  '
Set-TextCell 'D21' '17.01'
Set-TextCell 'E21' '  -3.35%  '
Set-TextCell 'D22' '5.977'
Set-TextCell 'E22' '  -6.19%  '
Set-TextCell 'D23' '27.994.63'
Set-TextCell 'E23' '  -1.79%  '
Set-TextCell 'D24' '11.25'
Set-TextCell 'E24' '  -4.79%  '
Set-TextCell 'D25' '2.137'
Set-TextCell 'E25' '  -1.64%  '
Set-TextCell 'D26' '157.67'
Set-TextCell 'E26' '  +2.25%  '
Set-TextCell 'D27' '20.14'
Set-TextCell 'E27' '  -4.43%  '
Set-TextCell 'D28' '1.986.42'
Set-TextCell 'E28' '  -2.23%  '
Set-TextCell 'D29' '2.130'
Set-TextCell 'E29' '  -10.70%  '
Set-TextCell 'D30' '124.59'
Set-TextCell 'E30' '  -4.42%  '
Set-TextCell 'D31' '1.155'
Set-TextCell 'E31' '  -6.17%  '
Set-TextCell 'D32' '5.629'
Set-TextCell 'E32' '  -4.41%  '
Set-TextCell 'D33' '0.08851'
Set-TextCell 'E33' '  -4.48%  '
Set-TextCell 'D34' '3.574'
Set-TextCell 'E34' '  -2.34%  '
Set-TextCell 'D35' '12.54'
Set-TextCell 'E35' '  -3.02%  '
Set-TextCell 'D36' '0.02310'
Set-TextCell 'E36' '  -2.25%  '
Set-TextCell 'D37' '0.2101'
Set-TextCell 'E37' '  -4.71%  '
Set-TextCell 'D40' '0.6353'
Set-TextCell 'E40' '  -4.67%  '
Set-TextCell 'D41' '1.187'
Set-TextCell 'E41' '  -0.61%  '
Set-TextCell 'D42' '0.9995'
Set-TextCell 'E42' '  +0.78%  '
Set-TextCell 'D43' '1.395'
Set-TextCell 'E43' '  -1.13%  '
Set-TextCell 'D44' '7.807'
Set-TextCell 'E44' '  -3.64%  '
Set-TextCell 'D45' '13.47'
Set-TextCell 'E45' '  -3.19%  '
Set-TextCell 'D46' '0.5906'
Set-TextCell 'E46' '  -3.73%  '
Set-TextCell 'D47' '3.685'
Set-TextCell 'E47' '  -1.83%  '
Set-TextCell 'D48' '2.004'
Set-TextCell 'E48' '  -1.94%  '
Set-TextCell 'D49' '122.26'
Set-TextCell 'E49' '  -4.39%  '
Set-TextCell 'D50' '1.186'
Set-TextCell 'E50' '  +2.98%  '
Set-TextCell 'D51' '0.06864'
Set-TextCell 'E51' '  -2.38%  '

Set-TextCell 'E17' '  -1.91%  '
Set-TextCell 'E20' '  +0.71%  '

# Rows 38 and 39 swapped rank position (Hedera <-> InternetComputer)
Set-TextCell 'B38' 'InternetComputer(DFINITY)'
Set-TextCell 'C38' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D38' '5.011'
Set-TextCell 'E38' '  -3.88%  '

Set-TextCell 'B39' 'Hedera'
Set-TextCell 'C39' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D39' '0.06030'
Set-TextCell 'E39' '  -3.55%  '
